$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 15 (the "checksum" field row under documents/file) - this shifts
# every subsequent row up by one, matching the target diff (dimension A1:N86 -> A1:N85).
$ws.Rows("15:15").Delete()
